# Update to CXF 3.1.0 (and refresh the "last modified" date footer field
# that PowerPoint stamps on every slide master / layout on save).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Footer "datetimeFigureOut" field: 15/05/30 -> 15/06/03
#    Present once on the slide master and once on every slide layout.
# ---------------------------------------------------------------------
$oldDate = "15/05/30"
$newDate = "15/06/03"

$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$layouts = $m.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $cl = $layouts.Item($li)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $sh = $cl.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                    $sh.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2. "Apache CXF 3.0.5" -> "Apache CXF 3.1.0" on slide 1.
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
for ($shi = 1; $shi -le $s1.Shapes.Count; $shi++) {
    $sh = $s1.Shapes.Item($shi)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            $t = $tr.Text
            $idx = $t.IndexOf("Apache CXF 3.0.5")
            if ($idx -ge 0) {
                $sub = $tr.Characters($idx + 1, 16)
                $sub.Text = "Apache CXF 3.1.0"
            }
        }
    }
}

# ---------------------------------------------------------------------
# 3. Merge the "Servlet " / "Container" runs into a single
#    "Servlet Container" run (slide 2).
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
for ($shi = 1; $shi -le $s2.Shapes.Count; $shi++) {
    $sh = $s2.Shapes.Item($shi)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "Servlet Container") {
                $sub = $tr.Characters(1, 17)
                $sub.Text = "Servlet Container"
            }
        }
    }
}
